# Reorder the outcome rows of the heterogeneity-analysis table (rows 3-12,
# columns A-G) into the new order requested by the commit:
#   new row 3  <- old row 4   (Male condom attitudes index)
#   new row 4  <- old row 3   (Used male condom at most recent sex (0/1))
#   new row 5  <- old row 10  (Contraceptive knowledge index)
#   new row 6  <- old row 11  (Modern contraceptive methods known (n))
#   new row 7  <- old row 12  (Discussed contraceptive use with recent partner (0/1))
#   new row 8  <- old row 5   (Can identify a female condom (0/1))
#   new row 9  <- old row 6   (Would be willing to try a female condom (0/1))
#   new row 10 <- old row 7   (Female condom attitudes index)
#   new row 11 <- old row 8   (Has ever used a female condom (0/1))
#   new row 12 <- old row 9   (Used a female condom in last 6 months (0/1))
#
# We use Copy / PasteSpecial (values) through a staging area far below the
# table so that entire rows (all of columns A:G) move together, preserving
# each cell's original text content (including leading/trailing spaces)
# and type without introducing any new number formats/styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

# Row mapping: destination row -> source row (in the original layout).
$rowMap = @{
    3  = 4
    4  = 3
    5  = 10
    6  = 11
    7  = 12
    8  = 5
    9  = 6
    10 = 7
    11 = 8
    12 = 9
}

$stagingStart = 100

# Step 1: stage a copy of each source row (A:G) below the table so the
# in-place permutation below never reads an already-overwritten cell.
# (The staging rows start out empty, but clear them first regardless so
# this step is self-contained / idempotent.)
foreach ($destRow in 3..12) {
    $srcRow = $rowMap[$destRow]
    $stageRow = $stagingStart + $srcRow
    $ws.Range("A$stageRow" + ":G$stageRow").ClearContents()
    $ws.Range("A$srcRow" + ":G$srcRow").Copy()
    $ws.Range("A$stageRow" + ":G$stageRow").PasteSpecial($xlPasteValues)
}

# Step 2: write each staged row into its final destination row. Clear the
# destination first: PasteSpecial(values) leaves existing destination
# content alone wherever the source cell is blank (e.g. the F/G "OR"
# columns on index rows), so a plain paste would leak stale values.
foreach ($destRow in 3..12) {
    $srcRow = $rowMap[$destRow]
    $stageRow = $stagingStart + $srcRow
    $ws.Range("A$destRow" + ":G$destRow").ClearContents()
    $ws.Range("A$stageRow" + ":G$stageRow").Copy()
    $ws.Range("A$destRow" + ":G$destRow").PasteSpecial($xlPasteValues)
}

# Step 3: clear the staging area.
$ws.Range("A$($stagingStart+3)" + ":G$($stagingStart+12)").Clear()

$excel.CutCopyMode = $false
